$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the NroPoliza / FechaSiniestro values for rows 2 and 3
$ws.Range("E2").Value = "'11111003115"
$ws.Range("G2").Value = "'30/04/2022"
$ws.Range("E3").Value = "'11111003131 "
$ws.Range("G3").Value = "'09/04/2021"

# Update column E width (target stored width 17.28515625 chars; the closest value
# reachable through the ColumnWidth COM property in this runtime is 17.33, reached
# by requesting 16.5)
$ws.Columns.Item(5).ColumnWidth = 16.5

# Update the selected cell/range in the sheet view
$ws.Range("E4").Select()
